$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Joy-Con analog stick gate (part link + name + qty)
$ws.Range("B5").Value = "https://konsolowo.pl/pl/czesci-naprawcze-joy-con/2133-analog-3d-galka-joystick-v2-joy-con-nintendo-switch-5903981901548.html"
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = "Gałka (muszkatołowa)"
$ws.Range("E5").Value = 2

# Row 6: DualSense adaptive triggers (hyperlinked part link + name + qty)
$ws.Range("C6").Value = "Triggery"
$ws.Range("E6").Value = 1
$ws.Range("B6").WrapText = $true
$null = $ws.Hyperlinks.Add($ws.Range("B6"), "https://konsolowo.pl/pl/dualsense-bdm-010/2810-triggery-adaptacyjne-haptyczne-l1-l2-r1-r2-dualsense-ps5-bdm-010-5903981916375.html")

# Row heights to match wrapped two-line content
$ws.Rows(5).RowHeight = 28.8
$ws.Rows(6).RowHeight = 28.8

# Match the saved selection/active cell
$null = $ws.Range("B6").Select()
